$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Cells.Item(40, 8).Value = 4266.3335
$ws.Cells.Item(40, 10).Value = 5999.5
$ws.Cells.Item(40, 12).Value = 5999.5
$ws.Cells.Item(40, 14).Value = -6349.5
# Row 64
$ws.Cells.Item(64, 8).Value = 4116.2856
$ws.Cells.Item(64, 9).Value = 0
$ws.Cells.Item(64, 10).Value = 4116.2856
$ws.Cells.Item(64, 11).Value = 0
$ws.Cells.Item(64, 12).Value = 4116.2856
$ws.Cells.Item(64, 13).ClearContents()
$ws.Cells.Item(64, 14).Value = -4612.2856
# Row 67
$ws.Cells.Item(67, 8).Value = 4116.2856
$ws.Cells.Item(67, 9).Value = 0
$ws.Cells.Item(67, 10).Value = 4116.2856
$ws.Cells.Item(67, 11).Value = 0
$ws.Cells.Item(67, 12).Value = 4116.2856
$ws.Cells.Item(67, 13).ClearContents()
$ws.Cells.Item(67, 14).Value = -5832.2856
# Row 74
$ws.Cells.Item(74, 8).Value = 8750
$ws.Cells.Item(74, 9).Value = 8000
$ws.Cells.Item(74, 11).Value = 8000
$ws.Cells.Item(74, 13).Value = -7064
# Row 77
$ws.Cells.Item(77, 8).Value = 8750
$ws.Cells.Item(77, 9).Value = 8000
$ws.Cells.Item(77, 11).Value = 40000
$ws.Cells.Item(77, 13).Value = -35320
# Row 111
$ws.Cells.Item(111, 8).Value = 525
$ws.Cells.Item(111, 9).Value = 525
$ws.Cells.Item(111, 10).Value = 0
$ws.Cells.Item(111, 11).Value = 1575
$ws.Cells.Item(111, 12).Value = 0
$ws.Cells.Item(111, 13).Value = 1492
$ws.Cells.Item(111, 14).ClearContents()

$ws = $wb.Worksheets.Item("ARM")
# Row 22
$ws.Cells.Item(22, 8).Value = 5387.5
$ws.Cells.Item(22, 9).Value = 516.6667
$ws.Cells.Item(22, 11).Value = 516.6667
$ws.Cells.Item(22, 13).Value = -217.6667
# Row 63
$ws.Cells.Item(63, 8).Value = 2333.3333
$ws.Cells.Item(63, 9).Value = 2333.3333
$ws.Cells.Item(63, 11).Value = 2333.3333
$ws.Cells.Item(63, 13).Value = -1647.3333
# Row 66
$ws.Cells.Item(66, 8).Value = 2333.3333
$ws.Cells.Item(66, 9).Value = 2333.3333
$ws.Cells.Item(66, 11).Value = 11666.6665
$ws.Cells.Item(66, 13).Value = -8234.666499999999
# Row 88
$ws.Cells.Item(88, 8).Value = 2034.6428
$ws.Cells.Item(88, 9).Value = 1079.6666
$ws.Cells.Item(88, 10).Value = 2750.875
$ws.Cells.Item(88, 11).Value = 1079.6666
$ws.Cells.Item(88, 12).Value = 2750.875
$ws.Cells.Item(88, 13).Value = -673.6666
$ws.Cells.Item(88, 14).Value = -3562.875
# Row 91
$ws.Cells.Item(91, 8).Value = 2034.6428
$ws.Cells.Item(91, 9).Value = 1079.6666
$ws.Cells.Item(91, 10).Value = 2750.875
$ws.Cells.Item(91, 11).Value = 1079.6666
$ws.Cells.Item(91, 12).Value = 2750.875
$ws.Cells.Item(91, 13).Value = 324.3334
$ws.Cells.Item(91, 14).Value = -5558.875
# Row 114
$ws.Cells.Item(114, 8).Value = 50000
$ws.Cells.Item(114, 10).Value = 50000
$ws.Cells.Item(114, 12).Value = 50000
$ws.Cells.Item(114, 14).Value = -58678
# Row 132
$ws.Cells.Item(132, 8).Value = 2744.875
$ws.Cells.Item(132, 10).Value = 400
$ws.Cells.Item(132, 12).Value = 1200
$ws.Cells.Item(132, 14).Value = -6260
# Row 134
$ws.Cells.Item(134, 8).Value = 56249.5
$ws.Cells.Item(134, 10).Value = 56249.5
$ws.Cells.Item(134, 12).Value = 56249.5
$ws.Cells.Item(134, 14).Value = -66389.5

$ws = $wb.Worksheets.Item("BSM")
# Row 15
$ws.Cells.Item(15, 8).Value = 0
$ws.Cells.Item(15, 9).Value = 0
$ws.Cells.Item(15, 10).Value = 0
$ws.Cells.Item(15, 11).Value = 0
$ws.Cells.Item(15, 12).Value = 0
$ws.Cells.Item(15, 13).ClearContents()
$ws.Cells.Item(15, 14).ClearContents()
# Row 35
$ws.Cells.Item(35, 8).Value = 14999
$ws.Cells.Item(35, 10).Value = 14999
$ws.Cells.Item(35, 12).Value = 14999
$ws.Cells.Item(35, 14).Value = -15619
# Row 86
$ws.Cells.Item(86, 8).Value = 4139.8
$ws.Cells.Item(86, 9).Value = 5179.6
$ws.Cells.Item(86, 10).Value = 3100
$ws.Cells.Item(86, 11).Value = 5179.6
$ws.Cells.Item(86, 12).Value = 3100
$ws.Cells.Item(86, 13).Value = -4056.6
$ws.Cells.Item(86, 14).Value = -5346
# Row 89
$ws.Cells.Item(89, 8).Value = 4139.8
$ws.Cells.Item(89, 9).Value = 5179.6
$ws.Cells.Item(89, 10).Value = 3100
$ws.Cells.Item(89, 11).Value = 25898
$ws.Cells.Item(89, 12).Value = 15500
$ws.Cells.Item(89, 13).Value = -20282
$ws.Cells.Item(89, 14).Value = -26732
# Row 110
$ws.Cells.Item(110, 8).Value = 107188
$ws.Cells.Item(110, 10).Value = 107188
$ws.Cells.Item(110, 12).Value = 107188
$ws.Cells.Item(110, 14).Value = -115368

$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Cells.Item(7, 8).Value = 222.14285
$ws.Cells.Item(7, 9).Value = 96.42856999999999
$ws.Cells.Item(7, 11).Value = 96.42856999999999
$ws.Cells.Item(7, 13).Value = 16.57143000000001
# Row 16
$ws.Cells.Item(16, 8).Value = 1682.25
$ws.Cells.Item(16, 9).Value = 1682.25
$ws.Cells.Item(16, 11).Value = 1682.25
$ws.Cells.Item(16, 13).Value = -1395.25
# Row 22
$ws.Cells.Item(22, 8).Value = 1800
$ws.Cells.Item(22, 9).Value = 2850
$ws.Cells.Item(22, 10).Value = 750
$ws.Cells.Item(22, 11).Value = 2850
$ws.Cells.Item(22, 12).Value = 750
$ws.Cells.Item(22, 13).Value = -2500
$ws.Cells.Item(22, 14).Value = -1450
# Row 62
$ws.Cells.Item(62, 8).Value = 5000
$ws.Cells.Item(62, 9).Value = 5000
$ws.Cells.Item(62, 11).Value = 5000
$ws.Cells.Item(62, 13).Value = -4376
# Row 65
$ws.Cells.Item(65, 8).Value = 5000
$ws.Cells.Item(65, 9).Value = 5000
$ws.Cells.Item(65, 11).Value = 25000
$ws.Cells.Item(65, 13).Value = -21880
# Row 113
$ws.Cells.Item(113, 8).Value = 1682.25
$ws.Cells.Item(113, 9).Value = 1682.25
$ws.Cells.Item(113, 11).Value = 1682.25
$ws.Cells.Item(113, 13).Value = 487.75

$ws = $wb.Worksheets.Item("CUL")
# Row 94
$ws.Cells.Item(94, 8).Value = 965.6667
$ws.Cells.Item(94, 9).Value = 965.6667
$ws.Cells.Item(94, 11).Value = 2897.0001
$ws.Cells.Item(94, 13).Value = -2221.0001
# Row 133
$ws.Cells.Item(133, 8).Value = 2000
$ws.Cells.Item(133, 9).Value = 2000
$ws.Cells.Item(133, 11).Value = 6000
$ws.Cells.Item(133, 13).Value = -940
# Row 134
$ws.Cells.Item(134, 8).Value = 125864.125
$ws.Cells.Item(134, 9).Value = 125864.125
$ws.Cells.Item(134, 11).Value = 377592.375
$ws.Cells.Item(134, 13).Value = -372522.375
# Row 136
$ws.Cells.Item(136, 8).Value = 2833.3333
$ws.Cells.Item(136, 9).Value = 3250
$ws.Cells.Item(136, 10).Value = 2000
$ws.Cells.Item(136, 11).Value = 9750
$ws.Cells.Item(136, 12).Value = 6000
$ws.Cells.Item(136, 13).Value = -4650
$ws.Cells.Item(136, 14).Value = -16200
# Row 137
$ws.Cells.Item(137, 8).Value = 479.8
$ws.Cells.Item(137, 9).Value = 479.8
$ws.Cells.Item(137, 11).Value = 1439.4
$ws.Cells.Item(137, 13).Value = 3660.6
# Row 139
$ws.Cells.Item(139, 8).Value = 2281.2856
$ws.Cells.Item(139, 9).Value = 795.2
$ws.Cells.Item(139, 11).Value = 2385.6
$ws.Cells.Item(139, 13).Value = 2754.4
# Row 140
$ws.Cells.Item(140, 8).Value = 2212
$ws.Cells.Item(140, 9).Value = 1111
$ws.Cells.Item(140, 11).Value = 3333
$ws.Cells.Item(140, 13).Value = 1847
# Row 141
$ws.Cells.Item(141, 8).Value = 2227.4
$ws.Cells.Item(141, 9).Value = 2227.4
$ws.Cells.Item(141, 11).Value = 6682.200000000001
$ws.Cells.Item(141, 13).Value = -1502.200000000001

$ws = $wb.Worksheets.Item("GSM")
# Row 53
$ws.Cells.Item(53, 8).Value = 8753
$ws.Cells.Item(53, 10).Value = 8753
$ws.Cells.Item(53, 12).Value = 8753
$ws.Cells.Item(53, 14).Value = -10015
# Row 70
$ws.Cells.Item(70, 8).Value = 56808224
$ws.Cells.Item(70, 9).Value = 68168670
$ws.Cells.Item(70, 10).Value = 5995
$ws.Cells.Item(70, 11).Value = 68168670
$ws.Cells.Item(70, 12).Value = 5995
$ws.Cells.Item(70, 13).Value = -68168400
$ws.Cells.Item(70, 14).Value = -6535
# Row 73
$ws.Cells.Item(73, 8).Value = 56808224
$ws.Cells.Item(73, 9).Value = 68168670
$ws.Cells.Item(73, 10).Value = 5995
$ws.Cells.Item(73, 11).Value = 68168670
$ws.Cells.Item(73, 12).Value = 5995
$ws.Cells.Item(73, 13).Value = -68167734
$ws.Cells.Item(73, 14).Value = -7867
# Row 113
$ws.Cells.Item(113, 8).Value = 1166.6666
$ws.Cells.Item(113, 9).Value = 1000
$ws.Cells.Item(113, 10).Value = 1250
$ws.Cells.Item(113, 11).Value = 1000
$ws.Cells.Item(113, 12).Value = 1250
$ws.Cells.Item(113, 13).Value = 1170
$ws.Cells.Item(113, 14).Value = -5590
# Row 126
$ws.Cells.Item(126, 8).Value = 9533.5
$ws.Cells.Item(126, 9).Value = 9533.5
$ws.Cells.Item(126, 10).Value = 0
$ws.Cells.Item(126, 11).Value = 28600.5
$ws.Cells.Item(126, 12).Value = 0
$ws.Cells.Item(126, 13).Value = -26130.5
$ws.Cells.Item(126, 14).ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Cells.Item(7, 8).Value = 14332.5
$ws.Cells.Item(7, 9).Value = 18748.75
$ws.Cells.Item(7, 10).Value = 5500
$ws.Cells.Item(7, 11).Value = 18748.75
$ws.Cells.Item(7, 12).Value = 5500
$ws.Cells.Item(7, 13).Value = -18636.75
$ws.Cells.Item(7, 14).Value = -5724
# Row 16
$ws.Cells.Item(16, 8).Value = 1750
$ws.Cells.Item(16, 9).Value = 1750
$ws.Cells.Item(16, 10).Value = 0
$ws.Cells.Item(16, 11).Value = 1750
$ws.Cells.Item(16, 12).Value = 0
$ws.Cells.Item(16, 13).Value = -1580
$ws.Cells.Item(16, 14).ClearContents()
# Row 22
$ws.Cells.Item(22, 8).Value = 5000
$ws.Cells.Item(22, 9).Value = 5500
$ws.Cells.Item(22, 10).Value = 3000
$ws.Cells.Item(22, 11).Value = 5500
$ws.Cells.Item(22, 12).Value = 3000
$ws.Cells.Item(22, 13).Value = -5205
$ws.Cells.Item(22, 14).Value = -3590
# Row 27
$ws.Cells.Item(27, 8).Value = 5000
$ws.Cells.Item(27, 9).Value = 5500
$ws.Cells.Item(27, 10).Value = 3000
$ws.Cells.Item(27, 11).Value = 5500
$ws.Cells.Item(27, 12).Value = 3000
$ws.Cells.Item(27, 13).Value = -5393
$ws.Cells.Item(27, 14).Value = -3214
# Row 40
$ws.Cells.Item(40, 8).Value = 6999.2
$ws.Cells.Item(40, 9).Value = 4999.5
$ws.Cells.Item(40, 11).Value = 4999.5
$ws.Cells.Item(40, 13).Value = -4863.5
# Row 46
$ws.Cells.Item(46, 8).Value = 4359.7812
$ws.Cells.Item(46, 9).Value = 7493
$ws.Cells.Item(46, 10).Value = 2479.85
$ws.Cells.Item(46, 11).Value = 7493
$ws.Cells.Item(46, 12).Value = 2479.85
$ws.Cells.Item(46, 13).Value = -7305
$ws.Cells.Item(46, 14).Value = -2855.85
# Row 126
$ws.Cells.Item(126, 8).Value = 14332.5
$ws.Cells.Item(126, 9).Value = 18748.75
$ws.Cells.Item(126, 10).Value = 5500
$ws.Cells.Item(126, 11).Value = 56246.25
$ws.Cells.Item(126, 12).Value = 16500
$ws.Cells.Item(126, 13).Value = -53776.25
$ws.Cells.Item(126, 14).Value = -21440
# Row 136
$ws.Cells.Item(136, 8).Value = 3002.0454
$ws.Cells.Item(136, 10).Value = 3323
$ws.Cells.Item(136, 12).Value = 9969
$ws.Cells.Item(136, 14).Value = -15069

$ws = $wb.Worksheets.Item("WVR")
# Row 15
$ws.Cells.Item(15, 8).Value = 5007498.5
$ws.Cells.Item(15, 10).Value = 14998
$ws.Cells.Item(15, 12).Value = 14998
$ws.Cells.Item(15, 14).Value = -15574
# Row 39
$ws.Cells.Item(39, 8).Value = 34666.668
$ws.Cells.Item(39, 10).Value = 34666.668
$ws.Cells.Item(39, 12).Value = 34666.668
$ws.Cells.Item(39, 14).Value = -35492.668
